$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Nakul", 1882.4, 3.5),
    @("Victor", 1644.9, 4),
    @("Sofia", 1634.3, 2),
    @("Sebastian", 1620, 2),
    @("Ahmed", 1604.4, 3),
    @("Francesco", 1595.4, 2),
    @("David", 1584, 1.5),
    @("Peter", 1502.1, 4.5),
    @("Zixiong", 1493.3, 3.5),
    @("Jaewon", 1470.3, 1.5),
    @("Niketan", 1460.1, 1.5),
    @("Zheng", 1459.8, 2),
    @("Jean-Pierre", 1453.7, 2.5),
    @("Gergo", 1448.8, 2.5),
    @("Patricia", 1411.3, 3.5),
    @("Nic", 1408.6, 3.5),
    @("Ruitao", 1399.7, 2),
    @("Hamdan", 1384.9, 3),
    @("Abdulrahman", 1378.5, 3),
    @("Theo", 1368.8, 2.5),
    @("Kirill", 1344.1, 3.5),
    @("Xiao", 1322.2, 1.5),
    @("Norah", 1310.7, 3.5),
    @("Steven", 1305.1, 2),
    @("Khalil", 1300.5, 2.5),
    @("Gabriel", 1258, 2.5),
    @("Aleksandra", 1252.8, 4),
    @("Abdulwahab", 1225.1, 2),
    @("Hesham", 1217.9, 2),
    @("Nour", 663.3, 1.5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
